$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: dates 41558 (2013-10-11) and 41559 (2013-10-12) with hours 2 and 4.5
$ws.Range("A11").Copy()
$ws.Range("A12:A13").PasteSpecial(-4122)

$ws.Range("A12").Value = 41558
$ws.Range("B12").Value = 2

$ws.Range("A13").Value = 41559
$ws.Range("B13").Value = 4.5

# Total formula row
$ws.Range("B28").Formula = "=SUM(B2:B27)"

# Update selection to mimic where the user left off
$ws.Range("B29").Select()
